# BooksWV.xlsx - "Ghost Writer HTML - OK"
#
# The underlying shared-strings table got re-ordered (9 bibliography blurbs
# for Institutos/Associações got their publisher name wrapped in markdown
# italics "*...*" and were moved within the string table), which is why the
# xml diff shows dozens of <v> index churn in column J (the "simples"
# column). Resolved against the shared strings, however, only 9 rows of
# column J actually change their visible text - each gains italics markup
# around the institution/publisher name. Setting the resolved cell values
# directly reproduces that outcome regardless of how the engine re-builds
# the shared string table on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J3").Value  = '**Vieira**, Waldo; ***700 Experimentos da Conscienciologia***; 1.058 p.; 40 seções; 100 subseções; 700 caps.; 28,5 x 21,5 x 7 cm; enc.; *Instituto Internacional de Projeciologia e Conscienciologia* (IIPC); Rio de Janeiro, RJ; 1994'
$ws.Range("J5").Value  = '**Vieira**, Waldo; ***Conscienciograma: Técnica de Avaliação da Consciência Integral***; 344 p.; 150 abrevs.; 106 assuntos das folhas de avaliação; 21 x 14 cm; br.; *Instituto Internacional de Projeciologia* (IIP); Rio de Janeiro, RJ; 1996'
$ws.Range("J7").Value  = '**Vieira**, Waldo; ***Minidefinições Conscienciais***; 164 p.; 450 minifrases; 15 x 10 cm; br.; *Instituto Internacional de Projeciologia* (IIP); Rio de Janeiro, RJ; 1996'
$ws.Range("J14").Value = '**Vieira**, Waldo; ***Homo sapiens pacificus***; 1.584 p.; 24 seções; 413 caps.; 29 x 21,5 x 7 cm; enc.; 3ª Ed. Gratuita; *Centro de Altos Estudos da Conscienciologia* (CEAEC); & *Associação Internacional Editares*; Foz do Iguaçu, PR; 2007'
$ws.Range("J16").Value = '**Vieira**, Waldo; ***Projeciologia: Panorama das Experiências da Consciência Fora do Corpo Humano***; 1.254 p.; 18 seções; 525 caps.; 28 x 21 x 7 cm; enc.; 10a Ed. rev. e aum.; *Associação Internacional Editares*; Foz do Iguaçu, PR; 2009'
$ws.Range("J18").Value = '**Vieira**, Waldo; ***Manual da Proéxis: Programação Existencial***; 164 p.; 40 caps.; 21 x 14 cm; br.; 5ª Ed. Ver.; *Associação Internacional Editares*; Foz do Iguaçu, PR; 2011'
$ws.Range("J20").Value = '**Vieira**, Waldo; ***Manual da Dupla Evolutiva***; 208 p.; 40 caps.; 21 x 14 cm; br.; 3ª Ed.; *Associação Internacional Editares*; Foz do Iguaçu, PR; 2012'
$ws.Range("J22").Value = '**Vieira**, Waldo; ***700 Experimentos da Conscienciologia***;  1.088 p.; 40 seções; 100 subseções; 700 caps..; 28,5 x 21,5 x 7 cm; enc.; *Associação Internacional Editares*; Foz do Iguaçu, PR; 2013'
$ws.Range("J24").Value = '**Vieira**, Waldo; ***Dicionário de Argumentos da Conscienciologia***; 1.572 p.;  651 caps.; 28,5 x 21,5 x 7 cm; enc.; *Associação Internacional Editares*; Foz do Iguaçu, PR; 2014'

# Active cell / selection moved from J9 to J6 in the saved view state.
$ws.Range("J6").Select()
